$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (pushes the existing "Total" row from 9 -> 10)
$ws.Rows.Item(9).Insert()

# The numeric-looking columns (dataPoints/BU/experimentos) are stored as text
# in this report, so force a text format before writing the new values to
# keep them as text rather than letting them be coerced into numbers.
# (Row 7 / SRSTR / ferrugem is untouched by this edit, so it is left as-is.)
$ws.Range("C2:E6").NumberFormat = "@"
$ws.Range("C8:E9").NumberFormat = "@"
$ws.Range("C10:E10").NumberFormat = "@"

# Update the data rows (2-6) with the new figures
$ws.Range("C2").Value = "80526"
$ws.Range("D2").Value = "336"
$ws.Range("E2").Value = "689"

$ws.Range("C3").Value = "11146"
$ws.Range("D3").Value = "68"
$ws.Range("E3").Value = "123"

$ws.Range("C4").Value = "26228"
$ws.Range("D4").Value = "154"
$ws.Range("E4").Value = "348"

$ws.Range("C5").Value = "35957"
$ws.Range("D5").Value = "200"
$ws.Range("E5").Value = "369"

$ws.Range("C6").Value = "11774"
$ws.Range("D6").Value = "59"
$ws.Range("E6").Value = "114"

# Row 7 (SRSTR / ferrugem) is unchanged

$ws.Range("C8").Value = "6818"
$ws.Range("D8").Value = "79"
$ws.Range("E8").Value = "132"

# New row 9: PRMDN / plantasMortas
$ws.Range("A9").Value = "PRMDN"
$ws.Range("B9").Value = "plantasMortas"
$ws.Range("C9").Value = "7740"
$ws.Range("D9").Value = "32"
$ws.Range("E9").Value = "63"

# Updated Total row, now at row 10
$ws.Range("A10").Value = "Total"
$ws.Range("B10").Value = "Total"
$ws.Range("C10").Value = "182828"
$ws.Range("D10").Value = "963"
$ws.Range("E10").Value = "1876"
